$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "GS-GM2": unhide helper columns B:D, then append new rows
# 258-266 (drum-kit category master list).
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("GS-GM2")
$ws2.Range("B1:D1").EntireColumn.Hidden = $false

$ws2.Range("A258").Value = 222
$ws2.Range("B258").Value = "Drum Kit Standard"
$ws2.Range("E258").Value = 256
$ws2.Range("F258").Value = "Drum Kit Standard"

$ws2.Range("A259").Value = 223
$ws2.Range("B259").Value = "Drum Kit Room"
$ws2.Range("E259").Value = 257
$ws2.Range("F259").Value = "Drum Kit Room"

$ws2.Range("A260").Value = 224
$ws2.Range("B260").Value = "Drum Kit Power"
$ws2.Range("E260").Value = 258
$ws2.Range("F260").Value = "Drum Kit Power"

$ws2.Range("A261").Value = 225
$ws2.Range("B261").Value = "Drum Kit Electronic"
$ws2.Range("E261").Value = 259
$ws2.Range("F261").Value = "Drum Kit Electronic"

$ws2.Range("A262").Value = 226
$ws2.Range("B262").Value = "Drum Kit Analog"
$ws2.Range("E262").Value = 260
$ws2.Range("F262").Value = "Drum Kit Analog"

$ws2.Range("A263").Value = 227
$ws2.Range("B263").Value = "Drum Kit Jazz"
$ws2.Range("E263").Value = 261
$ws2.Range("F263").Value = "Drum Kit Jazz"

$ws2.Range("A264").Value = 228
$ws2.Range("B264").Value = "Drum Kit Brush"
$ws2.Range("E264").Value = 262
$ws2.Range("F264").Value = "Drum Kit Brush"

$ws2.Range("A265").Value = 229
$ws2.Range("B265").Value = "Drum Kit Orchestra"
$ws2.Range("E265").Value = 263
$ws2.Range("F265").Value = "Drum Kit Orchestra"

$ws2.Range("A266").Value = 230
$ws2.Range("B266").Value = "Drum Kit SFX"
$ws2.Range("E266").Value = 264
$ws2.Range("F266").Value = "Drum Kit SFX"

# ------------------------------------------------------------------
# Sheet "GM2-XG": append new rows 482-530 (drum-kit detail mapping).
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("GM2-XG")

$ws3.Range("A482").Value = 256
$ws3.Range("B482").Value = "Drum Kit Standard"
$ws3.Range("G482").Value = 480
$ws3.Range("H482").Value = "Std1 Kit"

$ws3.Range("A483").Value = 256
$ws3.Range("G483").Value = 481
$ws3.Range("H483").Value = "Std2 Kit"

$ws3.Range("A484").Value = 256
$ws3.Range("G484").Value = 482
$ws3.Range("H484").Value = "Dry Kit"

$ws3.Range("A485").Value = 256
$ws3.Range("G485").Value = 483
$ws3.Range("H485").Value = "Bright Kit"

$ws3.Range("A486").Value = 256
$ws3.Range("G486").Value = 484
$ws3.Range("H486").Value = "Skim Kit"

$ws3.Range("A487").Value = 256
$ws3.Range("G487").Value = 485
$ws3.Range("H487").Value = "Slim Kit"

$ws3.Range("A488").Value = 256
$ws3.Range("G488").Value = 486
$ws3.Range("H488").Value = "Rogue Kit"

$ws3.Range("A489").Value = 256
$ws3.Range("G489").Value = 487
$ws3.Range("H489").Value = "Hob Kit"

$ws3.Range("A490").Value = 257
$ws3.Range("B490").Value = "Drum Kit Room"
$ws3.Range("G490").Value = 488
$ws3.Range("H490").Value = "Room Kit"

$ws3.Range("A491").Value = 257
$ws3.Range("G491").Value = 489
$ws3.Range("H491").Value = "Dark Kit"

$ws3.Range("A492").Value = 258
$ws3.Range("B492").Value = "Drum Kit Power"
$ws3.Range("G492").Value = 490
$ws3.Range("H492").Value = "Rock_Old Kit"

$ws3.Range("A493").Value = 258
$ws3.Range("G493").Value = 491
$ws3.Range("H493").Value = "Rock_Old2 Kit"

$ws3.Range("A494").Value = 259
$ws3.Range("B494").Value = "Drum Kit Electronic"
$ws3.Range("G494").Value = 492
$ws3.Range("H494").Value = "Electro Kit"

$ws3.Range("A495").Value = 260
$ws3.Range("B495").Value = "Drum Kit Analog"
$ws3.Range("G495").Value = 493
$ws3.Range("H495").Value = "Analog Kit"

$ws3.Range("A496").Value = 260
$ws3.Range("G496").Value = 494
$ws3.Range("H496").Value = "Analog2 Kit"

$ws3.Range("A497").Value = 260
$ws3.Range("G497").Value = 495
$ws3.Range("H497").Value = "Dance Kit"

$ws3.Range("A498").Value = 260
$ws3.Range("G498").Value = 496
$ws3.Range("H498").Value = "Hiphop Kit"

$ws3.Range("A499").Value = 260
$ws3.Range("G499").Value = 497
$ws3.Range("H499").Value = "Jungle Kit"

$ws3.Range("A500").Value = 260
$ws3.Range("G500").Value = 498
$ws3.Range("H500").Value = "Apogee Kit"

$ws3.Range("A501").Value = 260
$ws3.Range("G501").Value = 499
$ws3.Range("H501").Value = "Perigee Kit"

$ws3.Range("A502").Value = 261
$ws3.Range("B502").Value = "Drum Kit Jazz"
$ws3.Range("G502").Value = 500
$ws3.Range("H502").Value = "Jazz Kit"

$ws3.Range("A503").Value = 261
$ws3.Range("G503").Value = 501
$ws3.Range("H503").Value = "Jazz2 Kit"

$ws3.Range("A504").Value = 262
$ws3.Range("B504").Value = "Drum Kit Brush"
$ws3.Range("G504").Value = 502
$ws3.Range("H504").Value = "Brush Kit"

$ws3.Range("A505").Value = 262
$ws3.Range("G505").Value = 503
$ws3.Range("H505").Value = "Real_Brush Kit"

$ws3.Range("A506").Value = 263
$ws3.Range("B506").Value = "Drum Kit Orchestra"
$ws3.Range("G506").Value = 504
$ws3.Range("H506").Value = "Symphony Kit"

$ws3.Range("A507").Value = 259
$ws3.Range("B507").Value = "Drum Kit Electronic"
$ws3.Range("G507").Value = 505
$ws3.Range("H507").Value = "HipHop2 Kit"

$ws3.Range("A508").Value = 259
$ws3.Range("G508").Value = 506
$ws3.Range("H508").Value = "Break Kit"

$ws3.Range("A509").Value = 259
$ws3.Range("G509").Value = 507
$ws3.Range("H509").Value = "Tramp Kit"

$ws3.Range("A510").Value = 259
$ws3.Range("G510").Value = 508
$ws3.Range("H510").Value = "Amber Kit"

$ws3.Range("A511").Value = 259
$ws3.Range("G511").Value = 509
$ws3.Range("H511").Value = "Coffin Kit"

$ws3.Range("A512").Value = 256
$ws3.Range("B512").Value = "Drum Kit Standard"
$ws3.Range("G512").Value = 510
$ws3.Range("H512").Value = "Live_Std Kit"

$ws3.Range("A513").Value = 256
$ws3.Range("G513").Value = 511
$ws3.Range("H513").Value = "Live_Funk Kit"

$ws3.Range("A514").Value = 262
$ws3.Range("B514").Value = "Drum Kit Brush"
$ws3.Range("G514").Value = 512
$ws3.Range("H514").Value = "Live_Brush Kit"

$ws3.Range("A515").Value = 256
$ws3.Range("B515").Value = "Drum Kit Standard"
$ws3.Range("G515").Value = 513
$ws3.Range("H515").Value = "Live_Std_Perc Kit"

$ws3.Range("A516").Value = 256
$ws3.Range("G516").Value = 514
$ws3.Range("H516").Value = "Live_Funk_Perc Kit"

$ws3.Range("A517").Value = 262
$ws3.Range("B517").Value = "Drum Kit Brush"
$ws3.Range("G517").Value = 515
$ws3.Range("H517").Value = "Live_Brush_Perc Kit"

$ws3.Range("A518").Value = 264
$ws3.Range("B518").Value = "Drum Kit SFX"
$ws3.Range("G518").Value = 516
$ws3.Range("H518").Value = "SFX1 Kit"

$ws3.Range("A519").Value = 264
$ws3.Range("G519").Value = 517
$ws3.Range("H519").Value = "SFX2 Kit"

$ws3.Range("A520").Value = 259
$ws3.Range("B520").Value = "Drum Kit Electronic"
$ws3.Range("G520").Value = 518
$ws3.Range("H520").Value = "Techno_KS Kit"

$ws3.Range("A521").Value = 259
$ws3.Range("G521").Value = 519
$ws3.Range("H521").Value = "Techno_HI Kit"

$ws3.Range("A522").Value = 259
$ws3.Range("G522").Value = 520
$ws3.Range("H522").Value = "Techno_LO Kit"

$ws3.Range("A523").Value = 256
$ws3.Range("B523").Value = "Drum Kit Standard"
$ws3.Range("G523").Value = 521
$ws3.Range("H523").Value = "Sakura Kit"

$ws3.Range("A524").Value = 256
$ws3.Range("B524").Value = "Drum Kit Standard"
$ws3.Range("G524").Value = 522
$ws3.Range("H524").Value = "Small_Latin Kit"

$ws3.Range("A525").Value = 256
$ws3.Range("B525").Value = "Drum Kit Standard"
$ws3.Range("G525").Value = 523
$ws3.Range("H525").Value = "China Kit"

$ws3.Range("A526").Value = 256
$ws3.Range("B526").Value = "Drum Kit Standard"
$ws3.Range("G526").Value = 524
$ws3.Range("H526").Value = "Cuban Kit"

$ws3.Range("A527").Value = 256
$ws3.Range("G527").Value = 525
$ws3.Range("H527").Value = "Cuban2 Kit"

$ws3.Range("A528").Value = 256
$ws3.Range("G528").Value = 526
$ws3.Range("H528").Value = "Brazilian Kit"

$ws3.Range("A529").Value = 256
$ws3.Range("G529").Value = 527
$ws3.Range("H529").Value = "PopLatin1 Kit"

$ws3.Range("A530").Value = 256
$ws3.Range("G530").Value = 528
$ws3.Range("H530").Value = "PopLatin2 Kit"

# ------------------------------------------------------------------
# Fills: A482:B511 green (00B050), B512:B530 yellow (FFFF00).
# ------------------------------------------------------------------
$ws3.Range("B512:B530").Interior.Color = 65535
$ws3.Range("A482:B511").Interior.Color = 5287936

# ------------------------------------------------------------------
# Selection / active-sheet state, applied last-sheet-wins so the
# GM2-XG tab ends up active (matches workbook.xml activeTab).
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws1.Range("E25").Select()
$ws2.Range("E258:F266").Select()
$ws3.Range("H512").Select()
